$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cell "Save" in H1, matching the formatting of the existing
# header row (copy format from G1, which carries the bold/border/center style)
$ws.Range("H1").Value = "Save"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)

# Fill in the new "Save" column values for rows 2-7
$values = @(0, 1, 1, 0, 0, 1)
for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 8).Value = $values[$i]
}
